$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers
$ws.Range("B1").Value = "etapa 2"
$ws.Range("C1").Value = "etapa 3"
$ws.Range("D1").Value = "dcfvg"

# Row 2
$ws.Range("A2").Value = "sabor"
$ws.Range("B2").Value = 3
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = 3

# Row 3
$ws.Range("A3").Value = "precio"
$ws.Range("B3").Value = 20
$ws.Range("C3").Value = 3
$ws.Range("E3").Value = 30

# Row 4
$ws.Range("A4").Value = "sdfg"
